# Loan RBI, Variable Instalments
# Insert a new (blank) column before column M ("In Advance") on the
# "Repayment schedule" sheet, update a couple of values that moved as a
# result of the new Variable Instalment column, and leave the
# "Repayment schedule" tab as the active tab (it picks up tabSelected
# automatically from Summary when we select it here).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Select()

# Insert a blank column before the old "In Advance" column (M), shifting
# M..P to N..Q.
$ws.Range("M1").EntireColumn.Insert()

# K3 (Due) used to carry the full "interest + principal" figure; it now
# only reflects the interest-only portion for this variable-instalment
# period.
$ws.Range("K3").Value = 101.92

# N3 (old M3, "In Advance") now carries the new Variable Instalment
# principal amount instead of 0.
$ws.Range("N3").Value = 10000

$ws.Range("K8").Select() | Out-Null
